$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace A2 value and fill B2:J2 with "Unassigned"
$ws.Range("A2:J2").Value = "Unassigned"

# Row 3: A3 becomes "Microstomus kitt", clear rest of row (B3:J3 already empty)
$ws.Range("A3").Value = "Microstomus kitt"

# Delete old rows 4, 5, 6 (entire rows) - shifting cells up
$ws.Range("A4:J6").EntireRow.Delete()
